$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.954.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.589.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.482"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.245"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.590.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.946.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0719"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "198.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +8.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -8.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("E37").Value = "  +7.52%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.724.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -17.92%  "
